# Apply "added all test cases" edit:
#  - remove the old Sheet2 (its only value now lives in Sheet1 row 2 data)
#  - overwrite the Sheet1 test-data row (row 2) with the new submitted values
#  - widen column C to fit the new (longer) email address

$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

# Sheet2 only held a single throwaway value ("u") that is no longer needed.
$sheet2 = $wb.Worksheets.Item("Sheet2")
if ($sheet2) {
    $sheet2.Delete() | Out-Null
}

$ws = $wb.Worksheets.Item("Sheet1")

# New form-submission test data replacing the previous row.
$ws.Range("A2").Value = "bkjh"
$ws.Range("B2").Value = "jbkjhkj"
$ws.Range("C2").Value = "ghjhkhkjhkhjhj56@gmail.com"
$ws.Range("D2").Value = "bghjgkj@123"
$ws.Range("E2").Value = "bghjgkj@123"

# Column C needs to be a bit wider to comfortably fit the new email address.
$ws.Columns.Item(3).ColumnWidth = 27.08
